# The document has three inline "logo" pictures living in the page
# headers/footers (a Pearson logo duplicated in the default + first-page
# footers, and a BTEC logo in the first-page header). Their internal
# <wp:docPr>/name (and matching <pic:cNvPr>/name) attribute is out of
# sync with reality and needs to be swapped:
#   Pearson logo pictures : image1.png -> image2.png
#   BTEC logo picture     : image2.jpg -> image1.jpg
#
# InlineShape has no settable .Name property (true of real Word too), so
# rename by round-tripping through the floating Shape object, which does
# expose .Name, then convert back to an inline shape in place.
function Rename-InlineShape($inlineShape, [string]$newName) {
    $floatingShape = $inlineShape.ConvertToShape()
    $floatingShape.Name = $newName
    $floatingShape.ConvertToInlineShape() | Out-Null
}

$d = $word.ActiveDocument

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections($s)

    for ($h = 1; $h -le $section.Headers.Count; $h++) {
        $header = $section.Headers($h)
        if ($header.Exists) {
            $shapes = $header.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $pic = $shapes.Item($i)
                if ($pic.AlternativeText -eq "BTec_Logo-Orange") {
                    Rename-InlineShape $pic "image1.jpg"
                }
            }
        }
    }

    for ($f = 1; $f -le $section.Footers.Count; $f++) {
        $footer = $section.Footers($f)
        if ($footer.Exists) {
            $shapes = $footer.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $pic = $shapes.Item($i)
                if ($pic.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    Rename-InlineShape $pic "image2.png"
                }
            }
        }
    }
}
